$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (YearLevel), shifting YearLevel to column E
$ws.Columns("D").Insert()

# Set new column width to match the others
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# Fill in header and sample data for the new "Unit" column
$ws.Range("D1").Value = "Unit"
$ws.Range("D1").Font.Bold = $true

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "3"
$ws.Range("D2").Style = "Normal"
